$wb = $excel.ActiveWorkbook

# Update "展览" sheet (row 2: 126 -> 128, row 3: 46 -> 48 in column F)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 128
$wsExhibit.Range("F3").Value = 48

# Update "全部类型" sheet with the same values (duplicate data set)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 128
$wsAll.Range("F3").Value = 48
